# Update "想去人数" (number of people interested) figures in the
# 展览 (Exhibitions) and 全部类型 (All types) sheets, reflecting the
# freshly scraped data output.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 20886
    "F3"  = 811
    "F7"  = 7805
    "F9"  = 23
    "F12" = 53
    "F14" = 150
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}

# Rows differ slightly between the two sheets for the last three updates
# because "全部类型" has extra rows interleaved.
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F33").Value = 4973
$ws1.Range("F38").Value = 12944
$ws1.Range("F44").Value = 409

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F35").Value = 4973
$ws4.Range("F40").Value = 12944
$ws4.Range("F46").Value = 409
